# Update ObjTables metadata strings: bump version 0.0.8 -> 0.0.9,
# update the embedded date, and rename the `id=` attribute to `class=`
# on the per-table Data metadata rows.

$wb = $excel.ActiveWorkbook

$newDate = "2020-04-27 01:05:01"
$newVersion = "0.0.9"

# Sheet "!!_Table of contents": A1 (top banner) and A2 (TableOfContents metadata)
$wsToc = $wb.Worksheets.Item("!!_Table of contents")
$wsToc.Unprotect()
$wsToc.Range("A1").Value = "!!!ObjTables objTablesVersion='$newVersion' date='$newDate'"
$wsToc.Range("A2").Value = "!!ObjTables type='TableOfContents' tableFormat='row' description='Table of contents' date='$newDate' objTablesVersion='$newVersion'"
$wsToc.Protect()

# Sheet "!!_Schema": A1 Schema metadata
$wsSchema = $wb.Worksheets.Item("!!_Schema")
$wsSchema.Unprotect()
$wsSchema.Range("A1").Value = "!!ObjTables type='Schema' tableFormat='row' description='Table/model and column/attribute definitions' date='$newDate' objTablesVersion='$newVersion'"
$wsSchema.Protect()

# Sheet "!!Company": A1 Data metadata, id='Company' -> class='Company'
$wsCompany = $wb.Worksheets.Item("!!Company")
$wsCompany.Unprotect()
$wsCompany.Range("A1").Value = "!!ObjTables type='Data' tableFormat='column' class='Company' name='Companies' date='$newDate' objTablesVersion='$newVersion'"
$wsCompany.Protect()

# Sheet "!!People": A1 Data metadata, id='Person' -> class='Person'
$wsPeople = $wb.Worksheets.Item("!!People")
$wsPeople.Unprotect()
$wsPeople.Range("A1").Value = "!!ObjTables type='Data' tableFormat='row' class='Person' name='People' date='$newDate' objTablesVersion='$newVersion'"
$wsPeople.Protect()
